# Auto-generated edit script: updates crypto price/volume data in Sheet1
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.396.71"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.789.52"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.78"
$ws.Range("E5").Value = "  -2.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.553"
$ws.Range("E6").Value = "  -4.30%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.88"
$ws.Range("E8").Value = "  +4.32%  "
$ws.Range("E9").Value = "  -2.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0663"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0931"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "2.047.97"
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.04"
$ws.Range("E13").Value = "  +6.97%  "
$ws.Range("D14").Value = "1.795.23"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.640"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").Value = "34.419.40"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.25"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.18"
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "255.07"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("D20").Value = "0.0₃0747"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.41"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.23"
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  -4.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.94"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.06"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("E28").Value = "  -4.08%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.76"
$ws.Range("E30").Value = "  -4.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0513"
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.58"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.90"
$ws.Range("E34").Value = "  +4.71%  "
$ws.Range("D35").Value = "1.453.73"
$ws.Range("E35").Value = "  -5.15%  "
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0190"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.631"
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "83.22"
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.892"
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0507"
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.05"
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.89"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Value = "1.949.06"
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.30"
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "99.25"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.99"
$ws.Range("E51").Value = "  -3.62%  "
